$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Portuguese -> English genre translation, applied one term at a time
# (mirrors a sequence of whole-sheet Find & Replace passes) so that the
# shared-string table grows in the same order as the authored edit.

# Biografia -> Biography
foreach ($addr in @("F2")) {
    $ws.Range($addr).Value = 'Biography'
}

# Ação -> Action
foreach ($addr in @("F7","F8","F12","F15","F16","F17","F18","F20","F21","F23","F24","F25","F48","F49","F50","F51","F56","F57","F58","F78","F79","F84","F85","F86","F87","F88","F98","F99","F100","F101","F103","F104","F107","F108","F111","F115","F116","F117","F118","F120","F121","F124","F125","F134","F140","F141","F147","F148","F149","F150","F151","F155","F156","F157","F158","F173","F174","F181","F182","F183","F184","F186","F187","F188","F190","F191","F197","F198","F199","F200","F203","F206","F207","F210","F214","F215","F216","F217","F218","F219","F220","F222","F223","F224","F232")) {
    $ws.Range($addr).Value = 'Action'
}

# Animação -> Animation
foreach ($addr in @("F9","F10","F29","F30","F31","F32","F33","F34","F40","F41","F42","F43","F60","F61","F62","F63","F64","F65","F66","F67","F68","F69","F70","F71","F72","F93","F109","F110","F112","F126","F127","F129","F142","F143","F144","F145","F146","F160","F161","F162","F163","F164","F175","F176","F194","F195","F196","F208","F209","F233")) {
    $ws.Range($addr).Value = 'Animation'
}

# Suspense -> Thriller
foreach ($addr in @("F77","F123","F130","F131","F132","F165","F221","F231")) {
    $ws.Range($addr).Value = 'Thriller '
}

# Comédia -> Comedy
foreach ($addr in @("F102","F185")) {
    $ws.Range($addr).Value = 'Comedy'
}

# Romance -> Romantic
foreach ($addr in @("F26","F97","F152","F153","F192","F193")) {
    $ws.Range($addr).Value = 'Romantic'
}

# Ficcão Cientifica -> Science fiction
foreach ($addr in @("F28","F80","F94","F95","F96","F128","F225","F226","F227")) {
    $ws.Range($addr).Value = 'Science fiction'
}

# Terror -> Horror
foreach ($addr in @("F27","F52","F119")) {
    $ws.Range($addr).Value = 'Horror'
}

# Row 189 (Captain Marvel / #1 movie week) genre was corrected by hand to
# "Biography" instead of the bulk "Action" translation, which also reset
# the cells formatting to the default style.
$ws.Range("F189").ClearFormats()
$ws.Range("F189").Value = "Biography"

# Re-enable the AutoFilter over the full data range
$ws.Range("A1:J233").AutoFilter() | Out-Null
